$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the peer review text in C21 (Elizabeth Brooke Doyle's review of
#    Brandon Stehling / Matthew Hendrick): their scores changed from 1->2
#    and 2->3 respectively.
$c21Text = @"
Brandon Stehling: 
2
didn't do much but was there to try and help when asked.
Matthew Hendrick
3
did the map but that's it. seemed on top of things but didn't actually do much. 
"@
$ws.Range("C21").Value = $c21Text

# 2. Add the missing peer review text for group 9 (Lino Virgen / Andrew
#    Kennedy / Alexander Tesfazgi) into C26, which was previously blank.
$c26Text = @"
Lino Virgen
Rating: 5/5
Review:
Lino was our team leader and did his fair share of the work. He was great to work with, and he held me accountable to my share of the work.
Alexander Tesfazgi
Rating: 4.5/5
Review:
Alexander did his fair share of the work, but I didn't see him much during our scheduled class time. I do know he was in communication with Lino outside of class, but it would have been nice to discuss the project more with all three of us there.
"@
$ws.Range("C26").WrapText = $true
$ws.Range("C26").Value = $c26Text
$ws.Rows.Item(26).RowHeight = 165

# 3. Update the view state to match where the editor ended up scrolling /
#    selecting after making the edits.
$ws.Range("C28").Select()
